$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregar datos para psicologos en el proyecto - nueva fila 9 (Pablo)
$ws.Range("A9").Value = "Pablo"
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B9").Value = (Get-Date -Year 2025 -Month 5 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = "Primeros logros"
$ws.Range("E9").Value = "Estar conmigo"
$ws.Range("F9").Value = "En proceso"

$ws.Range("A10").Select()
